$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B13 was stored as an inline string "1"; convert it to a real number 1
$ws.Range("B13").Value = 1

# Append a new row of annotation data as row 14
$ws.Range("A14").Value = "Sunsi Wu"
$ws.Range("B14").Value = "4"
$ws.Range("C14").Value = "suggest"
$ws.Range("D14").Value = "SUG"
$ws.Range("E14").Value = "THE"
$ws.Range("F14").Value = "3a6bf25f-9f71-48b7-a40b-7e968e5f9337"
$ws.Range("G14").Value = "ry-TW-WAb_annotated.xlsx"
$ws.Range("H14").Value = "I suggest to change it to e.g. 'from the true to the approximate posterior' to avoid confusion."
